$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 06:17:59"
$c = $ws.Range("H2")
$c.Formula = '="96%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("N2").Value = "-2.7 °C 5:37 TU"
$ws.Range("O2").Value = "-1.7 °C"
$ws.Range("E3").Value = "2026-02-06 06:18:01"
$c = $ws.Range("H3")
$c.Formula = '="73%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("M3").Value = "-1.3 °C 5:35 TU"
$ws.Range("O3").Value = "-2.3 °C"
$ws.Range("E4").Value = "2026-02-06 06:18:04"
$c = $ws.Range("H4")
$c.Formula = '="62%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J4").Value = "993.6 hPa"
$ws.Range("N4").Value = "7.2 °C 5:51 TU"
$ws.Range("O4").Value = "12.0 °C"
$ws.Range("E5").Value = "2026-02-06 06:18:07"
$c = $ws.Range("H5")
$c.Formula = '="77%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J5").Value = "994.0 hPa"
$ws.Range("E6").Value = "2026-02-06 06:18:10"
$ws.Range("J6").Value = "995.1 hPa"
$ws.Range("N6").Value = "13.4 °C 5:49 TU"
$ws.Range("O6").Value = "14.3 °C"
$ws.Range("E7").Value = "2026-02-06 06:18:12"
$ws.Range("J7").Value = "994.9 hPa"
$ws.Range("N7").Value = "8.8 °C 5:38 TU"
$ws.Range("O7").Value = "9.9 °C"
$ws.Range("E8").Value = "2026-02-06 06:18:15"
$ws.Range("N8").Value = "4.0 °C 5:49 TU"
$ws.Range("O8").Value = "5.7 °C"
$ws.Range("E9").Value = "2026-02-06 06:18:18"
$ws.Range("N9").Value = "0.0 °C 5:52 TU"
$ws.Range("O9").Value = "2.0 °C"
$ws.Range("E10").Value = "2026-02-06 06:18:20"
$ws.Range("O10").Value = "4.9 °C"
$ws.Range("E11").Value = "2026-02-06 06:18:22"
$c = $ws.Range("H11")
$c.Formula = '="84%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J11").Value = "995.8 hPa"
$ws.Range("N11").Value = "1.2 °C 5:39 TU"
$ws.Range("O11").Value = "4.2 °C"
$ws.Range("E12").Value = "2026-02-06 06:18:25"
$ws.Range("O12").Value = "12.1 °C"
$ws.Range("E13").Value = "2026-02-06 06:18:28"
$ws.Range("N13").Value = "2.9 °C 5:50 TU"
$ws.Range("O13").Value = "6.0 °C"
$ws.Range("E14").Value = "2026-02-06 06:18:30"
$ws.Range("E15").Value = "2026-02-06 06:18:33"
$c = $ws.Range("H15")
$c.Formula = '="87%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J15").Value = "994.1 hPa"
$ws.Range("N15").Value = "2.4 °C 5:53 TU"
$ws.Range("O15").Value = "6.2 °C"
$ws.Range("E16").Value = "2026-02-06 06:18:36"
$ws.Range("E17").Value = "2026-02-06 06:18:39"
$ws.Range("J17").Value = "997.1 hPa"
$ws.Range("N17").Value = "1.0 °C 5:30 TU"
$ws.Range("O17").Value = "2.9 °C"
$ws.Range("E18").Value = "2026-02-06 06:18:42"
$ws.Range("L18").Value = "38.5 km/h - 291º 5:33 TU"
$ws.Range("E19").Value = "2026-02-06 06:18:44"
$c = $ws.Range("H19")
$c.Formula = '="97%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J19").Value = "997.4 hPa"
$ws.Range("E20").Value = "2026-02-06 06:18:47"
$c = $ws.Range("H20")
$c.Formula = '="70%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("O20").Value = "-2.3 °C"
$ws.Range("E21").Value = "2026-02-06 06:18:50"
$c = $ws.Range("H21")
$c.Formula = '="87%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J21").Value = "995.0 hPa"
$ws.Range("N21").Value = "2.1 °C 5:35 TU"
$ws.Range("O21").Value = "4.6 °C"
$ws.Range("E22").Value = "2026-02-06 06:18:53"
$c = $ws.Range("H22")
$c.Formula = '="85%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("N22").Value = "3.8 °C 5:58 TU"
$ws.Range("O22").Value = "7.4 °C"
$ws.Range("E23").Value = "2026-02-06 06:18:55"
$ws.Range("J23").Value = "994.2 hPa"
$ws.Range("N23").Value = "6.1 °C 5:59 TU"
$ws.Range("E24").Value = "2026-02-06 06:18:58"
$ws.Range("J24").Value = "993.1 hPa"
$ws.Range("O24").Value = "12.1 °C"
$ws.Range("E25").Value = "2026-02-06 06:19:00"
$ws.Range("J25").Value = "996.3 hPa"
$ws.Range("E26").Value = "2026-02-06 06:19:03"
$c = $ws.Range("H26")
$c.Formula = '="81%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("N26").Value = "-4.4 °C 5:58 TU"
$ws.Range("O26").Value = "-0.8 °C"
$ws.Range("E27").Value = "2026-02-06 06:19:06"
$ws.Range("J27").Value = "993.9 hPa"
$ws.Range("O27").Value = "7.2 °C"
$ws.Range("E28").Value = "2026-02-06 06:19:09"
$ws.Range("J28").Value = "997.1 hPa"
$ws.Range("N28").Value = "-0.7 °C 5:59 TU"
$ws.Range("O28").Value = "2.4 °C"
$ws.Range("E29").Value = "2026-02-06 06:19:11"
$c = $ws.Range("H29")
$c.Formula = '="66%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("N29").Value = "5.7 °C 5:46 TU"
$ws.Range("O29").Value = "10.8 °C"
$ws.Range("E30").Value = "2026-02-06 06:19:14"
$c = $ws.Range("H30")
$c.Formula = '="73%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("O30").Value = "-3.6 °C"
$ws.Range("E31").Value = "2026-02-06 06:19:16"
$ws.Range("J31").Value = "997.0 hPa"
$ws.Range("E32").Value = "2026-02-06 06:19:19"
$ws.Range("J32").Value = "995.5 hPa"
$ws.Range("N32").Value = "11.5 °C 5:44 TU"
$ws.Range("O32").Value = "14.6 °C"
$ws.Range("E33").Value = "2026-02-06 06:19:22"
$ws.Range("N33").Value = "4.3 °C 5:59 TU"
$ws.Range("O33").Value = "6.4 °C"
$ws.Range("E34").Value = "2026-02-06 06:19:24"
$c = $ws.Range("H34")
$c.Formula = '="82%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("O34").Value = "7.2 °C"
$ws.Range("E35").Value = "2026-02-06 06:19:27"
$ws.Range("N35").Value = "-3.4 °C 5:56 TU"
$ws.Range("E36").Value = "2026-02-06 06:19:29"
$c = $ws.Range("H36")
$c.Formula = '="68%"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("J36").Value = "996.9 hPa"
$ws.Range("N36").Value = "8.2 °C 5:50 TU"
$ws.Range("O36").Value = "11.2 °C"
$excel.CutCopyMode = $false
